# Swap the data (A, B, D, E, F, G, H, Q, R) between row 13 and row 14
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr13 = "$col" + "13"
    $addr14 = "$col" + "14"

    $val13 = $ws.Range($addr13).Value2
    $val14 = $ws.Range($addr14).Value2

    $ws.Range($addr13).Value2 = $val14
    $ws.Range($addr14).Value2 = $val13
}
